$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D (Price) and Column E (Volume(1h)) text updates, keyed by row number.
# $null in either slot means "leave unchanged". D values are written with a
# leading apostrophe where needed so Excel keeps them as literal text instead
# of re-parsing them as numbers (which would also silently eat trailing
# zeros, e.g. "5.330" -> 5.33). Values that can never parse as a plain number
# (multiple "." separators, e.g. "30.229.35") are left unprefixed since a
# prefix would needlessly mark the cell quotePrefix in the saved style.
$changes = @{
    2  = @("30.229.35", '  +0.31%  ')
    3  = @("1.860.52", '  -0.08%  ')
    4  = @($null, '  +0.01%  ')
    5  = @("'236.57", '  +1.13%  ')
    6  = @($null, '  +0.02%  ')
    7  = @("'0.4679", $null)
    8  = @("'0.2864", '  +1.36%  ')
    9  = @("'0.06540", '  +0.16%  ')
    10 = @("'21.94", '  +8.67%  ')
    11 = @("'0.07916", '  +1.16%  ')
    12 = @("'97.31", '  +1.64%  ')
    13 = @("1.867.61", $null)
    14 = @("'5.168", '  +0.94%  ')
    15 = @("'0.6833", '  +1.99%  ')
    16 = @("'270.69", '  -2.86%  ')
    17 = @("30.219.92", '  +0.19%  ')
    18 = @("'13.55", '  +7.23%  ')
    19 = @($null, '  +0.09%  ')
    20 = @("'0.000007348", '  +1.58%  ')
    21 = @("2.112.92", '  +0.67%  ')
    22 = @("'5.330", '  -2.15%  ')
    23 = @($null, '  +0.00%  ')
    24 = @("'6.172", '  +0.48%  ')
    25 = @("'167.67", '  +1.51%  ')
    26 = @("'9.224", '  -0.93%  ')
    27 = @($null, '  +0.06%  ')
    28 = @("'1.964", '  +3.30%  ')
    29 = @("'1.383", '  +2.82%  ')
    30 = @("'0.09864", '  +3.07%  ')
    31 = @("'4.366", '  -0.63%  ')
    32 = @($null, '  +0.49%  ')
    33 = @("'4.059", '  -1.30%  ')
    34 = @("'0.04707", '  +1.23%  ')
    35 = @("'1.133", '  +3.36%  ')
    36 = @("'0.7013", '  +0.20%  ')
    37 = @("'2.710", '  +0.20%  ')
    38 = @("'0.01871", '  +0.66%  ')
    39 = @($null, '  +3.83%  ')
    40 = @("'75.65", '  +4.59%  ')
    41 = @("'6.259", '  -0.34%  ')
    42 = @("'1.946", '  +1.06%  ')
    43 = @("'0.8510", '  -0.16%  ')
    46 = @("'103.42", '  +0.35%  ')
    47 = @("'956.89", '  -3.04%  ')
    48 = @("'7.176", '  +0.71%  ')
    49 = @("'9.220", '  +0.96%  ')
    50 = @("'34.15", '  +0.54%  ')
    51 = @("'0.05652", '  +0.54%  ')
}

foreach ($row in $changes.Keys) {
    $vals = $changes[$row]
    $dVal = $vals[0]
    $eVal = $vals[1]
    if ($null -ne $dVal) {
        $ws.Cells.Item($row, 4).Value = $dVal
    }
    if ($null -ne $eVal) {
        $ws.Cells.Item($row, 5).Value = $eVal
    }
}

# Rows 44 and 45 swap content (coin identity swaps along with updated values).
$ws.Cells.Item(44, 2).Value = 'TheSandbox'
$ws.Cells.Item(44, 3).Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
$ws.Cells.Item(44, 4).Value = "'0.4167"
$ws.Cells.Item(44, 5).Value = '  +0.27%  '

$ws.Cells.Item(45, 2).Value = 'PaxDollar'
$ws.Cells.Item(45, 3).Value = 'https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp'
$ws.Cells.Item(45, 4).Value = "'1.000"
$ws.Cells.Item(45, 5).Value = '  -0.01%  '
